# Colocando header nos graficos
$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlPasteAll = -4104

function Add-HeaderCell($ws, $addr, $text, $formatSourceAddr) {
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $src = $ws.Range($formatSourceAddr)
    $src.Copy()
    $cell.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
}

function Clear-LabelStyle($ws, $addr, $newText) {
    $cell = $ws.Range($addr)
    if ($newText -ne $null) {
        $cell.Value = $newText
    }
    $cell.ClearFormats()
}

# ---- Sheets 1-4: "Potencia Acumulada", "Geracao Periodo Medio",
#      "Atendimento a Ponta", "Potencia Incremental" share the same
#      row/column layout (A1 header + A2:A12 technology labels). ----
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Add-HeaderCell $ws "A1" "Fonte/Tecnologia" "B1"

    Clear-LabelStyle $ws "A2" $null
    Clear-LabelStyle $ws "A3" "Gás Natural"
    Clear-LabelStyle $ws "A4" "Carvão"
    Clear-LabelStyle $ws "A5" $null
    Clear-LabelStyle $ws "A6" "Óleos Comb"
    Clear-LabelStyle $ws "A7" $null
    Clear-LabelStyle $ws "A8" "Eólica"
    Clear-LabelStyle $ws "A9" $null
    Clear-LabelStyle $ws "A10" $null
    Clear-LabelStyle $ws "A11" "Pot. Compl."
    Clear-LabelStyle $ws "A12" $null
}

# ---- Sheet 5: "Emissoes Totais (MtCO2eq)" ----
$ws5 = $wb.Worksheets.Item(5)

Add-HeaderCell $ws5 "A1" "Período" "B1"

Clear-LabelStyle $ws5 "A2" "P.Médio"
Clear-LabelStyle $ws5 "A3" "P.Crítico"

# Remove the "Teto" row entirely (row 4)
$ws5.Rows.Item(4).Delete()

# ---- Sheet 6: "Custo Total (bilhões de R$)" ----
$ws6 = $wb.Worksheets.Item(6)

Add-HeaderCell $ws6 "A1" "Tipo Expansão" "A2"

# B1 header text changes from "Custo" to "2015" but must stay text,
# not be auto-converted to a number by Excel's smart typing.
$b1 = $ws6.Range("B1")
$b1.NumberFormat = "@"
$b1.Value = "2015"
$ws6.Range("A2").Copy()
$b1.PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Clear-LabelStyle $ws6 "A2" "Expansão Centralizada"
$ws6.Range("B2").Value = 579

Clear-LabelStyle $ws6 "A3" "Expansão por GD"
$ws6.Range("B3").Value = 99
